$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 (year 2025) metrics with the refreshed BIBI data
$ws.Range("C8").Value = 1027
$ws.Range("D8").Value = 168
$ws.Range("E8").Value = 859
$ws.Range("F8").Value = 6.890894175553733
$ws.Range("G8").Value = 83.64167478091528
$ws.Range("H8").Value = 16.35832521908471
